$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D40: replace the old "Ivan: preguntar..." note with the ireport link
$ws.Range("D40").Value = "http://mygnet.net/articulos/java/creacion_de_graficos_en_ireport.707"

# Row 52 "Estetica - botones imagenes" gains Responsable + 50% progress
$ws.Range("B52").Value = "Agustina"
$ws.Range("C52").Value = 0.5
$ws.Range("C52").NumberFormat = "0%"

# Row 53 "Estetica - etiquetas..." gains Responsable + 100% progress
$ws.Range("B53").Value = "Agustina"
$ws.Range("C53").Value = 1
$ws.Range("C53").NumberFormat = "0%"

# Row 56 "Reunion con Ivan y Josefina" marked complete
$ws.Range("C56").Value = 1
$ws.Range("C56").NumberFormat = "0%"

# Row 58 "Agregar los comandos abajo..." moves from "en proceso" text to 100% done
$ws.Range("C58").Value = 1
$ws.Range("C58").NumberFormat = "0%"

# Row 61 "No asignar cliente en venta de factura B" marked complete
$ws.Range("C61").Value = 1
$ws.Range("C61").NumberFormat = "0%"

# New tasks appended at the bottom of the list
$ws.Range("A62").Value = "Error en iva cuando consumidor final en ticket"
$ws.Range("B62").Value = "Lucas"

$ws.Range("A63").Value = "Error en listado de control - pagos de la fecha"
$ws.Range("B63").Value = "Agustina"

$ws.Range("A64").Value = "Migracion de datos"
$ws.Range("B64").Value = "Agustina"

$ws.Range("A65").Value = "Paginado de consultas"

# Row 66 previously held the orphaned "Ivan: preguntar..." note far below the table;
# it now holds a new task instead
$ws.Range("A66").Value = "Reporte listado de cliente"

$ws.Range("A67").Value = "Reporte para contador"
$ws.Range("B67").Value = "Agustina"

# Leave the cursor where the author left it, ready for the next entry
$ws.Range("B68").Select() | Out-Null
